# Nota Contabila - clear the "Concedii medicale din fonduri" suma value and
# remove the now-unneeded TOTAL row for the "Fonduri handicap" section
# ("zile suspendare - needs checking").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Suma value for "Concedii medicale din fonduri" (row 15, col F)
$ws.Range("F15").ClearContents()

# Delete the trailing TOTAL row (row 40) which summed the "Fonduri handicap"
# section (row 39); this row is dropped entirely in the revised template.
$ws.Rows(40).Delete()

# Restore the selection to match the updated view
$ws.Range("L29").Select()
